# This script applies a re-shuffle of the weekly price records on Sheet1.
# Rows 2, 5 and 6 rotate their data among themselves, and rows 4 and 8 swap
# their data. Columns A-C, F-K are identical across the affected rows, so
# only D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), R (Origen), S (Precio $/Kg) and
# T (Kg / unidad) need to be rewritten per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now holds what used to be row 6's data.
$ws.Range("D2").Value = 44174
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("Q2").Value = "$/bandeja 18 kilos"
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1083
$ws.Range("T2").Value = 18

# Row 4 now holds what used to be row 8's data.
$ws.Range("D4").Value = 44160
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("Q4").Value = "$/bandeja 18 kilos"
$ws.Range("R4").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S4").Value = 1361
$ws.Range("T4").Value = 18

# Row 5 now holds what used to be row 2's data.
$ws.Range("D5").Value = 44169
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21000
$ws.Range("Q5").Value = "$/bandeja 18 kilos"
$ws.Range("R5").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S5").Value = 1167
$ws.Range("T5").Value = 18

# Row 6 now holds what used to be row 5's data.
$ws.Range("D6").Value = 44524
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 27000
$ws.Range("O6").Value = 28000
$ws.Range("P6").Value = 27500
$ws.Range("Q6").Value = "$/bandeja 18 kilos"
$ws.Range("R6").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S6").Value = 1528
$ws.Range("T6").Value = 18

# Row 8 now holds what used to be row 4's data.
$ws.Range("D8").Value = 44533
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 140
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = "$/caja 10 kilos"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1450
$ws.Range("T8").Value = 10
